$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.661.02'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '3.400.30'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.09'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '3.400.88'
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.565'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.119'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.427'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.80%  '
$ws.Range('D13').Value = '3.982.38'
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.02'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000173'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.00%  '
$ws.Range('D17').Value = '63.730.06'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').Value = '3.380.65'
$ws.Range('E18').Value = '  -3.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.11'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('E22').Value = '  -3.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.517'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000115'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.03%  '
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.99'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('E31').Value = '  -5.57%  '
$ws.Range('E32').Value = '  -2.94%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '22.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.94'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.25%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.53'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.46%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.33'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.62%  '
$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.818'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.16%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.98%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0727'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.78%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.781.17'
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.64'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.37'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.32%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.40'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.45%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.39'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.20%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0303'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.32%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.36'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +10.83%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '322.46'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.15%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.03'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.10%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.32'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.33%  '
